$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(-156.374, -1.848, 0.065, -322.229, 9.481999999999999)
    3  = @(-118.148, -6.341, 0, -154.666, -81.63)
    4  = @(-69.48, -3.728, 0, -106.008, -32.953)
    5  = @(-3.458, -0.68, 0.497, -13.426, 6.51)
    6  = @(7.267, 6.489, 0, 5.072, 9.461)
    7  = @(8.551, 7.634, 0, 6.355, 10.746)
    8  = @(26.43, 4.805, 0, 15.65, 37.211)
    9  = @(-5.798, -4.321, 0, -8.428000000000001, -3.168)
    10 = @(-11.218, -8.791, 0, -13.719, -8.717000000000001)
    11 = @(-3.87, -0.54, 0.589, -17.916, 10.176)
    12 = @(11.44, 6.58, 0, 8.032, 14.848)
    13 = @(13.388, 8.083, 0, 10.142, 16.634)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    $ws.Cells.Item($row, 4).Value = $rowValues[0]
    $ws.Cells.Item($row, 5).Value = $rowValues[1]
    $ws.Cells.Item($row, 6).Value = $rowValues[2]
    $ws.Cells.Item($row, 7).Value = $rowValues[3]
    $ws.Cells.Item($row, 8).Value = $rowValues[4]
}
